# Edit script: restructure "Estado_actual_flujo" sheet (remove fecha_fin/DEFAULT
# column, merge into a single "fecha" column), and change the active sheet
# selection back to "Estado_actual_flujo".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 = "Estado_actual_flujo"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estado_actual_flujo")

# Clear out the whole used range first so no stale cells (old layout) remain.
$ws2.Range("A1:L8").Clear()

# --- Header row ---
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "fecha"
$ws2.Range("C1").Value = "activo"
$ws2.Range("D1").Value = "id_registro"
$ws2.Range("E1").Value = "id_estado"

$ws2.Range("A1:E1").Font.Bold = $true

$insertPrefix = "INSERT INTO sgr.estado_actual_flujo(id_estado_actual, fecha, activo, id_registro, id_estado) values ("

# --- Row 2 ---
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "2017-09-09 14:29:08.308049"
$ws2.Range("C2").Value = "'TRUE"
$ws2.Range("D2").Value = 11
$ws2.Range("E2").Value = 2
$ws2.Range("G2").Value = $insertPrefix
$ws2.Range("H2").Formula = '=A2&","&"''"&B2&"''"&"::TIMESTAMP"&","&"''"&C2&"''"&","&D2&","&E2&")"'
$ws2.Range("J2").Formula = '=G2&H2&";"'
$ws2.Range("K2").Value = $ws2.Range("J2").Text
$ws2.Range("K2").Font.Color = 255

# --- Row 3 ---
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "2017-09-09 16:38:03.599437"
$ws2.Range("C3").Value = "'TRUE"
$ws2.Range("D3").Value = 12
$ws2.Range("E3").Value = 2
$ws2.Range("G3").Value = $insertPrefix
$ws2.Range("H3").Formula = '=A3&","&"''"&B3&"''"&"::TIMESTAMP"&","&"''"&C3&"''"&","&D3&","&E3&")"'
$ws2.Range("J3").Formula = '=G3&H3&";"'
$ws2.Range("K3").Value = $ws2.Range("J3").Text
$ws2.Range("K3").Font.Color = 255

# --- Row 4 ---
$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "2017-09-30 09:37:15.278727"
$ws2.Range("C4").Value = "'TRUE"
$ws2.Range("D4").Value = 14
$ws2.Range("E4").Value = 2
$ws2.Range("G4").Value = $insertPrefix
$ws2.Range("H4").Formula = '=A4&","&"''"&B4&"''"&"::TIMESTAMP"&","&"''"&C4&"''"&","&D4&","&E4&")"'
$ws2.Range("J4").Formula = '=G4&H4&";"'
$ws2.Range("K4").Value = $ws2.Range("J4").Text
$ws2.Range("K4").Font.Color = 255

# --- Trailing styled (empty) cells, rows 5-8 in column K ---
$ws2.Range("K5").Font.Color = 255
$ws2.Range("K6").Font.Color = 255
$ws2.Range("K7").Font.Color = 255
$ws2.Range("K8").Font.Color = 255

# --- Column widths matching the new layout ---
$ws2.Columns.Item(8).ColumnWidth = 57.14

# --- Selection / active sheet ---
$ws2.Activate()
$ws2.Range("K2:K4").Select()
